$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "latest week" block (rows 292-293),
# pushing the existing rows 292-316 down to 294-318.
$ws.Rows("292:293").Insert()

# Row 292: new weekly entry - "Primera" quality
$ws.Range("A292").Value = 3
$ws.Range("B292").Value = "Femacal de La Calera"
$ws.Range("C292").Value = "Coquimbo"
$ws.Range("D292").Value = 44714
$ws.Range("E292").Value = 5
$ws.Range("F292").Value = "Fruta"
$ws.Range("G292").Value = 100101
$ws.Range("H292").Value = "Berries"
$ws.Range("I292").Value = 100112025
$ws.Range("J292").Value = "Frutilla"
$ws.Range("K292").Value = "Sin especificar"
$ws.Range("L292").Value = "Primera"
$ws.Range("M292").Value = 40
$ws.Range("N292").Value = 7000
$ws.Range("O292").Value = 7000
$ws.Range("P292").Value = 7000
$ws.Range("Q292").Value = "$/bandeja 7 kilos"
$ws.Range("R292").Value = "Provincia de Melipilla"
$ws.Range("S292").Value = 1000
$ws.Range("T292").Value = 7

# Row 293: new weekly entry - "Segunda" quality
$ws.Range("A293").Value = 3
$ws.Range("B293").Value = "Femacal de La Calera"
$ws.Range("C293").Value = "Coquimbo"
$ws.Range("D293").Value = 44714
$ws.Range("E293").Value = 5
$ws.Range("F293").Value = "Fruta"
$ws.Range("G293").Value = 100101
$ws.Range("H293").Value = "Berries"
$ws.Range("I293").Value = 100112025
$ws.Range("J293").Value = "Frutilla"
$ws.Range("K293").Value = "Sin especificar"
$ws.Range("L293").Value = "Segunda"
$ws.Range("M293").Value = 45
$ws.Range("N293").Value = 5000
$ws.Range("O293").Value = 5000
$ws.Range("P293").Value = 5000
$ws.Range("Q293").Value = "$/bandeja 7 kilos"
$ws.Range("R293").Value = "Provincia de Melipilla"
$ws.Range("S293").Value = 714
$ws.Range("T293").Value = 7
